# refactor: web client and wifi scanning procedure.
# Updates to the "Typography" sheet:
#  - Rows 4-6 (Default / Large / Small typographies) gain/change their
#    "Wildcard Characters" (col G) and "Wildcard Ranges" (col I) values so
#    they match the values already used by the APs_SSID / modalWindowTitle
#    rows (wifi-scanning related widgets now reuse the same wildcard set).
#  - Rows 9-10 (Clock / modalWindowTitle) fonts (col C) are corrected from
#    "APs_RSSI" (a stray/incorrect font name) to "Default".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Typography")

$wildcardChars  = '~`!@#$%^&*()_-+={}[]:"|;' + [char]39 + '\<>?,./'
$wildcardRanges = '0-9,A-Z,a-z'

$ws.Range("G4").Value = $wildcardChars
$ws.Range("I4").Value = $wildcardRanges

$ws.Range("G5").Value = $wildcardChars
$ws.Range("I5").Value = $wildcardRanges

$ws.Range("G6").Value = $wildcardChars
$ws.Range("I6").Value = $wildcardRanges

$ws.Range("C9").Value = "Default"
$ws.Range("C10").Value = "Default"
